$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: B1 -> supplier_id, C1 -> stok_tanggal, D1 -> stok_jumlah
$ws.Range("B1").Value = "supplier_id"
$ws.Range("C1").Value = "stok_tanggal"
$ws.Range("D1").Value = "stok_jumlah"

# Update row 2, A2 value from 13 to 16
$ws.Range("A2").Value = 16

# Clear rows 3 and 4: B:D entirely (so the cells disappear), A column contents only
$ws.Range("B3:D4").Clear()
$ws.Range("A3:A4").ClearContents()

# Update selection
$ws.Range("B3:D4").Select()
